$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 (model selection results changed) ---

# Row 2: top model (full interaction model, now includes day_night, dropped ar1 habitat term)
$ws.Range("C2").Value = "mean_accel ~ habitat_type * season * day_night + (1 | animal_id) + ar1(season + 0 | animal_id)"
$ws.Range("D2").Value = "m"
$ws.Range("E2").Value = 18539
$ws.Range("F2").Value = 0.368452446462205
$ws.Range("G2").Value = 4002.95079718084
$ws.Range("H2").Value = -7869.90159436169
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = -7337.62262514587
$ws.Range("L2").Value = 18471

# Row 3: now season + day_night model (m7)
$ws.Range("C3").Value = "mean_accel ~ season + day_night + (1 | animal_id) + season:day_night"
$ws.Range("D3").Value = "m7"
$ws.Range("E3").Value = 18539
$ws.Range("F3").Value = 0.376870079516635
$ws.Range("G3").Value = 3729.37131394347
$ws.Range("H3").Value = -7422.74262788694
$ws.Range("I3").Value = 447.158966474748
$ws.Range("J3").Value = [double]"7.95543924150648e-98"
$ws.Range("K3").Value = -7281.84525368275
$ws.Range("L3").Value = 18521

# Row 4: now habitat_type + season model (m4)
$ws.Range("C4").Value = "mean_accel ~ habitat_type + season + (1 | animal_id) + habitat_type:season"
$ws.Range("D4").Value = "m4"
$ws.Range("E4").Value = 18539
$ws.Range("F4").Value = 0.383810591934266
$ws.Range("G4").Value = 3279.01963096275
$ws.Range("H4").Value = -6518.0392619255
$ws.Range("I4").Value = 1351.86233243619
$ws.Range("J4").Value = [double]"2.79784960010847e-294"
$ws.Range("K4").Value = -6361.48662392085
$ws.Range("L4").Value = 18519

# --- New rows 5-9 (additional candidate models) ---

# Row 5: season-only model (m2)
$ws.Range("A5").Value = "lognormal"
$ws.Range("B5").Value = "log"
$ws.Range("C5").Value = "mean_accel ~ season + (1 | animal_id)"
$ws.Range("D5").Value = "m2"
$ws.Range("E5").Value = 18539
$ws.Range("F5").Value = 0.387027149613895
$ws.Range("G5").Value = 3182.15379064369
$ws.Range("H5").Value = -6352.30758128738
$ws.Range("I5").Value = 1517.59401307431
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -6305.34178988598
$ws.Range("L5").Value = 18533

# Row 6: habitat_type + day_night model (m5)
$ws.Range("A6").Value = "lognormal"
$ws.Range("B6").Value = "log"
$ws.Range("C6").Value = "mean_accel ~ habitat_type + day_night + (1 | animal_id) + habitat_type:day_night"
$ws.Range("D6").Value = "m5"
$ws.Range("E6").Value = 18539
$ws.Range("F6").Value = 0.381863469281159
$ws.Range("G6").Value = 2460.89752126732
$ws.Range("H6").Value = -4877.79504253465
$ws.Range("I6").Value = 2992.10655182704
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -4705.58714072953
$ws.Range("L6").Value = 18517

# Row 7: habitat_type + day_night model (m6)
$ws.Range("A7").Value = "lognormal"
$ws.Range("B7").Value = "log"
$ws.Range("C7").Value = "mean_accel ~ habitat_type + day_night + (1 | animal_id) + habitat_type:day_night"
$ws.Range("D7").Value = "m6"
$ws.Range("E7").Value = 18539
$ws.Range("F7").Value = 0.381863469281159
$ws.Range("G7").Value = 2460.89752126732
$ws.Range("H7").Value = -4877.79504253465
$ws.Range("I7").Value = 2992.10655182704
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -4705.58714072953
$ws.Range("L7").Value = 18517

# Row 8: day_night-only model (m3)
$ws.Range("A8").Value = "lognormal"
$ws.Range("B8").Value = "log"
$ws.Range("C8").Value = "mean_accel ~ day_night + (1 | animal_id)"
$ws.Range("D8").Value = "m3"
$ws.Range("E8").Value = 18539
$ws.Range("F8").Value = 0.387667847897727
$ws.Range("G8").Value = 2286.96015148016
$ws.Range("H8").Value = -4561.92030296032
$ws.Range("I8").Value = 3307.98129140137
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -4514.95451155893
$ws.Range("L8").Value = 18533

# Row 9: habitat_type-only model (m1)
$ws.Range("A9").Value = "lognormal"
$ws.Range("B9").Value = "log"
$ws.Range("C9").Value = "mean_accel ~ habitat_type + (1 | animal_id)"
$ws.Range("D9").Value = "m1"
$ws.Range("E9").Value = 18539
$ws.Range("F9").Value = 0.390232291554728
$ws.Range("G9").Value = 1913.36153476976
$ws.Range("H9").Value = -3812.72306953952
$ws.Range("I9").Value = 4057.17852482217
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = -3757.92964623789
$ws.Range("L9").Value = 18532
